$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Img" header in F1, matching the header style used by B1:D1.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Img"

# Image hyperlinks for rows 2-8 (order matches the rIds 1-7 from the diff).
# Hyperlinks.Add seeds the cell text + the <hyperlink display="..."> with
# TextToDisplay; we set it to the raw URL first, then overwrite the cell's
# visible text with the friendly caption afterward so the stored `display`
# attribute keeps the URL while the cell shows the caption.
$ws.Hyperlinks.Add($ws.Range("F4"), "https://raw.githubusercontent.com/hvijay31/Diet/main/roti.png", "", "", "https://raw.githubusercontent.com/hvijay31/Diet/main/roti.png")
$ws.Range("F4").Value = "roti.png (501×498) (raw.githubusercontent.com)"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://raw.githubusercontent.com/hvijay31/Diet/main/rice.png", "", "", "https://raw.githubusercontent.com/hvijay31/Diet/main/rice.png")
$ws.Range("F3").Value = "rice.png (499×499) (raw.githubusercontent.com)"

$ws.Hyperlinks.Add($ws.Range("F2"), "https://raw.githubusercontent.com/hvijay31/Diet/main/brrice.png", "", "", "https://raw.githubusercontent.com/hvijay31/Diet/main/brrice.png")
$ws.Range("F2").Value = "brrice.png (494×505) (raw.githubusercontent.com)"

$ws.Hyperlinks.Add($ws.Range("F6"), "https://raw.githubusercontent.com/hvijay31/Diet/main/mutton.png", "", "", "https://raw.githubusercontent.com/hvijay31/Diet/main/mutton.png")
$ws.Range("F6").Value = "mutton.png (494×505) (raw.githubusercontent.com)"

$ws.Hyperlinks.Add($ws.Range("F7"), "https://raw.githubusercontent.com/hvijay31/Diet/main/chicken.png", "", "", "https://raw.githubusercontent.com/hvijay31/Diet/main/chicken.png")
$ws.Range("F7").Value = "chicken.png (494×505) (raw.githubusercontent.com)"

$ws.Hyperlinks.Add($ws.Range("F5"), "https://raw.githubusercontent.com/hvijay31/Diet/main/idili.png", "", "", "https://raw.githubusercontent.com/hvijay31/Diet/main/idili.png")
$ws.Range("F5").Value = "idili.png (494×505) (raw.githubusercontent.com)"

$ws.Hyperlinks.Add($ws.Range("F8"), "https://raw.githubusercontent.com/hvijay31/Diet/main/dosa.png", "", "", "https://raw.githubusercontent.com/hvijay31/Diet/main/dosa.png")
$ws.Range("F8").Value = "dosa.png (494×505) (raw.githubusercontent.com)"

# Match the cursor position captured in the saved workbook.
$ws.Range("F8").Select()
